# Add two new rows (104, 105) of daily data to each of the 8 worksheets.
# Row 104 -> date serial 45967 (2025-11-06)
# Row 105 -> date serial 45968 (2025-11-07), remn_amt = 0

$wb = $excel.ActiveWorkbook

# New B-column values per sheet, in sheet order (index 1..8), for row 104 and row 105.
$newValues = @{
    1 = @(10558528, 0)
    2 = @(14243367, 0)
    3 = @(3534238, 0)
    4 = @(996050, 0)
    5 = @(1616520, 0)
    6 = @(1917123, 0)
    7 = @(302352, 0)
    8 = @(302009, 0)
}

$dateRow104 = 45967
$dateRow105 = 45968

for ($i = 1; $i -le 8; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $vals = $newValues[$i]
    $b104 = $vals[0]
    $b105 = $vals[1]

    # copy the existing date-format from the previous date cell (A103)
    $dateFormat = $ws.Range("A103").NumberFormat

    $ws.Range("A104").Value = $dateRow104
    $ws.Range("A104").NumberFormat = $dateFormat
    $ws.Range("B104").Value = $b104

    $ws.Range("A105").Value = $dateRow105
    $ws.Range("A105").NumberFormat = $dateFormat
    $ws.Range("B105").Value = $b105
}
